# Populate the "Recipe Category" (column D) and "Food Category" (column E)
# values for every recipe row (rows 2-22) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @{ Row = 2;  D = "";                          E = "Veg" },
    @{ Row = 3;  D = "Breakfast, Lunch, Dinner";   E = "Vegetarian, Veg" },
    @{ Row = 4;  D = "";                          E = "" },
    @{ Row = 5;  D = "Breakfast";                  E = "" },
    @{ Row = 6;  D = "";                          E = "Vegan, Veg" },
    @{ Row = 7;  D = "";                          E = "Vegan, Veg" },
    @{ Row = 8;  D = "Breakfast, Lunch, Snack";    E = "Vegetarian, Veg" },
    @{ Row = 9;  D = "";                          E = "" },
    @{ Row = 10; D = "";                          E = "" },
    @{ Row = 11; D = "Breakfast";                  E = "" },
    @{ Row = 12; D = "";                          E = "" },
    @{ Row = 13; D = "Lunch";                      E = "" },
    @{ Row = 14; D = "Snack";                      E = "Veg" },
    @{ Row = 15; D = "Snack";                      E = "Vegetarian, Veg" },
    @{ Row = 16; D = "Dinner";                     E = "" },
    @{ Row = 17; D = "Breakfast";                  E = "" },
    @{ Row = 18; D = "Breakfast, Snack";           E = "" },
    @{ Row = 19; D = "Snack";                      E = "Vegan, Veg" },
    @{ Row = 20; D = "";                          E = "" },
    @{ Row = 21; D = "";                          E = "" },
    @{ Row = 22; D = "";                          E = "" }
)

foreach ($item in $data) {
    $dCell = $ws.Cells.Item($item.Row, 4)
    if ($item.D -eq "") {
        # A bare "'" forces the cell to a (zero-length) text value instead of
        # clearing it outright; ClearFormats() drops the quote-prefix style
        # that the leading apostrophe otherwise leaves behind.
        $dCell.Value = "'"
        $dCell.ClearFormats()
    } else {
        $dCell.Value = $item.D
    }

    $eCell = $ws.Cells.Item($item.Row, 5)
    if ($item.E -eq "") {
        $eCell.Value = "'"
        $eCell.ClearFormats()
    } else {
        $eCell.Value = $item.E
    }
}
